$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.929.53"
$ws.Range("E2").Value = "  -1.18%  "

$ws.Range("D3").Value = "1.993.82"
$ws.Range("E3").Value = "  -2.87%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.55"
$ws.Range("E5").Value = "  -2.68%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.603"
$ws.Range("E6").Value = "  -2.92%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.12"
$ws.Range("E8").Value = "  -5.20%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.375"
$ws.Range("E9").Value = "  -2.32%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0778"
$ws.Range("E10").Value = "  +0.78%  "

$ws.Range("E11").Value = "  -3.10%  "

$ws.Range("D12").Value = "2.293.31"
$ws.Range("E12").Value = "  -2.59%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.00"
$ws.Range("E13").Value = "  -4.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.02"
$ws.Range("E14").Value = "  -3.16%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.733"
$ws.Range("E15").Value = "  -3.13%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.07"
$ws.Range("E16").Value = "  -4.06%  "

$ws.Range("D17").Value = "1.991.13"
$ws.Range("E17").Value = "  -3.12%  "

$ws.Range("D18").Value = "36.827.61"
$ws.Range("E18").Value = "  -1.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.03"
$ws.Range("E19").Value = "  -0.91%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.37"
$ws.Range("E20").Value = "  -1.79%  "

$ws.Range("D21").Value = "0.0₃0810"
$ws.Range("E21").Value = "  -1.51%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "221.87"
$ws.Range("E22").Value = "  -1.98%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.19%  "

$ws.Range("E24").Value = "  +1.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.16"
$ws.Range("E25").Value = "  -6.97%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.72"
$ws.Range("E26").Value = "  -2.49%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.06"
$ws.Range("E27").Value = "  -8.43%  "

$ws.Range("E28").Value = "  -3.77%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.48"
$ws.Range("E29").Value = "  -3.82%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.33"
$ws.Range("E30").Value = "  -1.26%  "

$ws.Range("E31").Value = "  -4.14%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.45"
$ws.Range("E32").Value = "  -1.66%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0605"
$ws.Range("E33").Value = "  -2.77%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.38"
$ws.Range("E34").Value = "  -4.34%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.31"
$ws.Range("E35").Value = "  -6.82%  "

$ws.Range("E36").Value = "  +1.17%  "

$ws.Range("E37").Value = "  +0.04%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.11"
$ws.Range("E38").Value = "  -5.06%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.28"
$ws.Range("E39").Value = "  -0.70%  "

$ws.Range("D40").Value = "1.470.20"
$ws.Range("E40").Value = "  -0.51%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0214"
$ws.Range("E41").Value = "  -4.28%  "

$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0916"
$ws.Range("E42").Value = "  -4.00%  "

$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "93.76"
$ws.Range("E43").Value = "  -4.56%  "

$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.27"
$ws.Range("E44").Value = "  -2.08%  "

$ws.Range("E45").Value = "  -4.78%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.11"
$ws.Range("E46").Value = "  -6.35%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.14"
$ws.Range("E47").Value = "  -1.70%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.997"
$ws.Range("E48").Value = "  -3.04%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.89"
$ws.Range("E49").Value = "  -1.79%  "

$ws.Range("D50").Value = "2.185.85"
$ws.Range("E50").Value = "  -2.42%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.96"
$ws.Range("E51").Value = "  -4.12%  "
